$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data for "Enterprises (absolute #)" (row 12) and
# "Enterprises density (per 1000 people)" (row 13) was regenerated with
# the density figure listed first. Swap the label/value pairs between
# the two rows so the density row now precedes the absolute-count row,
# matching the regenerated shared-string order.
$labelA12 = $ws.Range("A12").Value2
$valueD12 = $ws.Range("D12").Value2
$fmtD12   = $ws.Range("D12").NumberFormat

$labelA13 = $ws.Range("A13").Value2
$valueD13 = $ws.Range("D13").Value2
$fmtD13   = $ws.Range("D13").NumberFormat

# Force the numeric-looking labels to be written back as text (they were
# stored as shared strings, not numbers, in the source file) by
# temporarily switching to a text number format while assigning them.
$ws.Range("A12").Value2 = $labelA13

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = $valueD13
$ws.Range("D12").NumberFormat = $fmtD12

$ws.Range("A13").Value2 = $labelA12

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = $valueD12
$ws.Range("D13").NumberFormat = $fmtD13
